$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column stays text (it contains values like "1.000" or
# "30.712.27" that Excel would otherwise reinterpret as numbers/dates).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.712.27"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "1.890.23"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "247.50"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4940"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "0.2956"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "0.06816"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "1.890.58"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "17.23"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "92.17"
$ws.Range("E12").Value = "  +6.30%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07248"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6786"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.069"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "30.695.04"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "0.000007962"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "2.135.98"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "4.838"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "191.40"
$ws.Range("E23").Value = "  +33.64%  "
$ws.Range("D24").Value = "6.064"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").Value = "9.397"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "156.87"
$ws.Range("E26").Value = "  +4.64%  "
$ws.Range("D27").Value = "19.04"
$ws.Range("E27").Value = "  +11.30%  "
$ws.Range("D28").Value = "1.917"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "1.404"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "4.311"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").Value = "0.08982"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "4.020"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "0.05187"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "0.7437"
$ws.Range("E34").Value = "  +3.84%  "
$ws.Range("D35").Value = "1.119"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").Value = "2.721"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").Value = "0.01843"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "2.681"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "2.155"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "0.9416"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "0.4437"
$ws.Range("E41").Value = "  +3.97%  "
$ws.Range("D42").Value = "105.70"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "5.756"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "7.663"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").Value = "0.1342"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("E47").Value = "  +3.15%  "
$ws.Range("D48").Value = "1.431"
$ws.Range("E48").Value = "  +6.87%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.3950"
$ws.Range("E49").Value = "  +3.75%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.645"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").Value = "33.52"
$ws.Range("E51").Value = "  +2.46%  "
